$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Update Sheet1 header row to use bracket unit notation instead of parenthesis notation.
$ws1.Range("A1").Value = "C1 [m/s]"
$ws1.Range("B1").Value = "C2 [mol/L]"
$ws1.Range("C1").Value = "C3 [(m^3)/kg]"

# Leave the active selection on C1 of Sheet1, matching the edited cell.
$ws1.Activate()
$ws1.Range("C1").Select()
